$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.758.31"
$ws.Range("E2").Value = "  +4.01%  "

$ws.Range("D3").Value = "'1.870.42"
$ws.Range("E3").Value = "  +2.88%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "'277.16"
$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").Value = "'0.5279"
$ws.Range("E7").Value = "  +3.55%  "

$ws.Range("D8").Value = "'0.3418"
$ws.Range("E8").Value = "  -3.25%  "

$ws.Range("D9").Value = "'0.06930"
$ws.Range("E9").Value = "  +3.91%  "

$ws.Range("D10").Value = "'20.00"
$ws.Range("E10").Value = "  -0.38%  "

$ws.Range("D11").Value = "'0.8023"
$ws.Range("E11").Value = "  -3.12%  "

$ws.Range("D12").Value = "'0.07766"
$ws.Range("E12").Value = "  -1.81%  "

$ws.Range("D13").Value = "'1.881.70"
$ws.Range("E13").Value = "  +4.31%  "

$ws.Range("D14").Value = "'90.21"
$ws.Range("E14").Value = "  +2.85%  "

$ws.Range("D15").Value = "'5.175"
$ws.Range("E15").Value = "  +1.93%  "

$ws.Range("D16").Value = "'14.55"
$ws.Range("E16").Value = "  +3.24%  "

$ws.Range("D17").Value = "'1.002"
$ws.Range("E17").Value = "  +0.20%  "

$ws.Range("D18").Value = "'0.000008029"
$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").Value = "'26.800.91"
$ws.Range("E20").Value = "  +4.01%  "

$ws.Range("D21").Value = "'2.093.93"
$ws.Range("E21").Value = "  +2.79%  "

$ws.Range("D22").Value = "'4.745"
$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("E23").Value = "  +0.34%  "

$ws.Range("D24").Value = "'6.159"
$ws.Range("E24").Value = "  +0.87%  "

$ws.Range("D25").Value = "'2.363"
$ws.Range("E25").Value = "  +6.26%  "

$ws.Range("D26").Value = "'146.35"
$ws.Range("E26").Value = "  +3.20%  "

$ws.Range("D27").Value = "'17.33"
$ws.Range("E27").Value = "  +1.32%  "

$ws.Range("D28").Value = "'1.656"
$ws.Range("E28").Value = "  -0.98%  "

$ws.Range("D29").Value = "'113.36"
$ws.Range("E29").Value = "  +3.71%  "

$ws.Range("D30").Value = "'4.335"
$ws.Range("E30").Value = "  +0.33%  "

$ws.Range("D31").Value = "'4.324"
$ws.Range("E31").Value = "  +1.99%  "

$ws.Range("D32").Value = "'0.08907"
$ws.Range("E32").Value = "  +1.56%  "

$ws.Range("D33").Value = "'0.04931"
$ws.Range("E33").Value = "  +0.92%  "

$ws.Range("D34").Value = "'1.165"
$ws.Range("E34").Value = "  +2.33%  "

$ws.Range("D35").Value = "'0.7282"
$ws.Range("E35").Value = "  -0.07%  "

$ws.Range("D36").Value = "'2.886"
$ws.Range("E36").Value = "  +0.88%  "

$ws.Range("D37").Value = "'3.269"
$ws.Range("E37").Value = "  +4.30%  "

$ws.Range("D38").Value = "'0.01852"
$ws.Range("E38").Value = "  -0.13%  "

$ws.Range("D39").Value = "'2.319"
$ws.Range("E39").Value = "  -2.51%  "

$ws.Range("D40").Value = "'0.5135"
$ws.Range("E40").Value = "  -0.58%  "

$ws.Range("D41").Value = "'0.9488"
$ws.Range("E41").Value = "  -1.85%  "

$ws.Range("D42").Value = "'115.98"
$ws.Range("E42").Value = "  +4.29%  "

$ws.Range("D43").Value = "'6.158"
$ws.Range("E43").Value = "  -1.00%  "

$ws.Range("D44").Value = "'8.082"
$ws.Range("E44").Value = "  +0.64%  "

$ws.Range("D45").Value = "'1.001"
$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("D46").Value = "'0.4460"
$ws.Range("E46").Value = "  -2.28%  "

$ws.Range("D47").Value = "'0.1339"
$ws.Range("E47").Value = "  -1.97%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.06051"
$ws.Range("E48").Value = "  +3.69%  "

$ws.Range("D49").Value = "'36.31"
$ws.Range("E49").Value = "  -0.52%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.259"
$ws.Range("E50").Value = "  +0.60%  "

$ws.Range("D51").Value = "'1.487"
$ws.Range("E51").Value = "  -0.94%  "
